$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = @(-0.013353625880095286,-0.014286745014970487,-0.01485427826917625,-0.011379887435091853,-0.024855316282409848,-0.019484112134156142,-0.01598694341085038,-0.005742209496297201,-0.024311529802082087,-0.022070568924394986,-0.0220440642326439,-0.01550586339668339,-0.010179593944384575,-0.006531812455454684,-0.02387168611486599,-0.03742635216340434,-0.014669638585753176,-0.02725966295578253,-0.00854289495105042,-0.025143217144756175,-0.016058742413290403,-0.010575639395727653,-0.01263906833536683,-0.0219868380041365,-0.019411081964487226,-0.02220358758902996,-0.037490637102067764,-0.012166891681991428,-0.03278645235662175,-0.009172672318530911,-0.017618637861345043,-0.013715625904831564,-0.028012293924302115,-0.020837035127614874,-0.02198372838623224,-0.015464180266425695,-0.016037769543771445,-0.020935652237031355,-0.02343987775984972,-0.017455712922869314,-0.03497695292579907,-0.013704868258657117,-0.009580378124288864,-0.015503753526335226,-0.004479924047585403,-0.014695008444309524,-0.01902842295576513,-0.013573337943315826,-0.017404679660836433,-0.01519506167533022,-0.014671631512403178,-0.02144625277336059,-0.02743977465035865,-0.008135683154992149,-0.017185391241795404,-0.012490569853519488,-0.016988862969480257,-0.009437600151676048,-0.005724717480569858,-0.004834123603580999,-0.011096971780559024,-0.018067340304572136,-0.013507198310030347,-0.010390656271193432,-0.015801740999593823,-0.013353670568753534,-0.0239416845498125,-0.02479157296654393,-0.015405383222005364,-0.011328035499101727,-0.017309727049769897,-0.02135728890941074,-0.016487600518937183,-0.008605134562180825,-0.029957942397707944,-0.0319637376743553,-0.008962915165400733,-0.011592704508772289,-0.008889209045737407,-0.016744038398397262,-0.02279170067519619,-0.027176772568400422,-0.019852767204630684,-0.02530152090110076,-0.010299210151462557,-0.012083276848252682,-0.011334917177616272,-0.01919983973538938,-0.010757942931194833,-0.02144902500733291,-0.020199556741492604,-0.023534170865208388,-0.017456703391012504,-0.011467438857035203,-0.02724122406871182,-0.011991064844808908,-0.011822787537477472,-0.012563671988061518,-0.024035749667920898,-0.01375477072040404)
$row2 = @(-0.0015242804047821406,-0.0015682410913208295,-0.0016273211972283078,-0.0012057615679017767,-0.0027887486606467675,-0.0021523222742965663,-0.0018023419060136001,-0.0006572524823470455,-0.002783000003216232,-0.0024262477490738124,-0.0024171382171028306,-0.001753773371596866,-0.0011640017877548327,-0.0006576621965135117,-0.002637816731246584,-0.004174003626491481,-0.001650850640944885,-0.0031173331937969346,-0.0009715956413658247,-0.002778740791389183,-0.0018185094791325616,-0.0011993044175407997,-0.0014455389053131577,-0.0024919283526912404,-0.002211207147456266,-0.0024418631578607326,-0.004192959813250488,-0.001307947450409614,-0.0037299966421302207,-0.0010253863593865772,-0.0019565680473586974,-0.0015155554610171373,-0.003178509080782331,-0.0023117505020573523,-0.0025018183661108637,-0.0016776663923837338,-0.0018317132687677012,-0.0023915960021250204,-0.0026800350289554005,-0.001915653398428831,-0.004003099729528377,-0.0015655759601631578,-0.0010962062695301954,-0.0017287407232506446,-0.0005115963327734113,-0.0016822037160373138,-0.002131993556378937,-0.0015331728492422138,-0.001977238170524854,-0.001728038199941541,-0.0016456176074703642,-0.002428071168110517,-0.003030809611018965,-0.0009068227828290574,-0.0019427055311881131,-0.0014183762159510542,-0.0019314103316654264,-0.0010214832808054988,-0.0006516740628853178,-0.000553378489174958,-0.0012673284400729564,-0.0020118469557894813,-0.001543368810352038,-0.00117983044887914,-0.0017795055440703329,-0.0015161421663180354,-0.002630064343989496,-0.0027950952804286023,-0.0017529157526240664,-0.0012260757676444615,-0.0019812241805309766,-0.002368068184686693,-0.0018176498405149011,-0.0009777370219087235,-0.0033383875491694204,-0.003586370275840321,-0.0010149002576209443,-0.0013215973496071278,-0.0009460507319081564,-0.0018150598157612434,-0.0025780607930181968,-0.003077197334253563,-0.0022040969345365497,-0.0028960750667359992,-0.0011746070837645869,-0.0013801222460955672,-0.0012804101091285721,-0.002133539768595715,-0.001230932827687436,-0.0024517144942276053,-0.002273374734914042,-0.002649985756690102,-0.0019978692128400518,-0.0012801971911260521,-0.003063559873396129,-0.0013074558859032508,-0.0013416632029810854,-0.0013831957122753355,-0.0027514470805758075,-0.0014915223272152134)

$data = New-Object "object[,]" 2,100
for ($i = 0; $i -lt 100; $i++) {
    $data[0, $i] = $row1[$i]
    $data[1, $i] = $row2[$i]
}

$ws.Range("A1:CV2").Value = $data
